# Fix the swapped hotel_info / review_info sheets:
#  - The physical first sheet currently holds "hotel_info" data but (per the
#    target) should become "review_info" (header row only, 25 columns).
#  - The physical second sheet currently holds "review_info" data but should
#    become "hotel_info" (header + 1 data row, with a new "State" column
#    inserted right after "Hotel_Name").
#
# We rename via a temporary name to avoid a duplicate-name collision, then
# rewrite each sheet's cell content to match the corrected layout.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename sheets (swap names) -------------------------------------------
$ws1.Name = "__tmp_swap__"
$ws2.Name = "hotel_info"
$ws1.Name = "review_info"

# --- Rebuild "review_info" sheet (physically the former sheet 1) ----------
$ws1.Cells.Clear()

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value2 = $reviewHeaders[$i]
}

# --- Rebuild "hotel_info" sheet (physically the former sheet 2) -----------
$ws2.Cells.Clear()

$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value2 = $hotelHeaders[$i]
}

$ws2.Cells.Item(2, 1).Value2 = 62034
$ws2.Cells.Item(2, 2).Value2 = "La Quinta Inns & Suites Boutte"
$ws2.Cells.Item(2, 3).Value2 = "Louisiana"
$ws2.Cells.Item(2, 4).Value2 = "Boutte"
$ws2.Cells.Item(2, 5).Value2 = 70039
$ws2.Cells.Item(2, 6).Value2 = "https://www.tripadvisor.com/Hotel_Review-g40046-d3475909-Reviews-La_Quinta_Inn_Suites_Boutte-Boutte_Louisiana.html"
$ws2.Cells.Item(2, 7).Value2 = "La Quinta Inn & Suites Boutte"

# English_Reviews_num ("176") and Total_Reviews_num ("181") are stored as
# text in the source data (not numbers), so force text formatting before
# assigning the numeric-looking strings - otherwise they'd be auto-coerced
# to numbers.
$ws2.Cells.Item(2, 8).NumberFormat = "@"
$ws2.Cells.Item(2, 8).Value2 = "176"
# Local_Rank (column I / 9) intentionally left blank
$ws2.Cells.Item(2, 10).NumberFormat = "@"
$ws2.Cells.Item(2, 10).Value2 = "181"
